$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.675700187683105
$ws.Range("B1").Value = 2.415036916732788
$ws.Range("C1").Value = 2.120635032653809
$ws.Range("D1").Value = 1.763802647590637
$ws.Range("E1").Value = 1.671298623085022
